$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume table (row swaps + updated quotes).
# A few D-column values are digit strings that Excel would otherwise auto-
# parse into a number (losing trailing zeros such as "5.430" -> 5.43, or the
# thousands-dot formatting). Force those cells to Text before writing the
# literal string, then drop back to the Normal style so no stray number
# formatting is left behind.

$ws.Range("D2").Value = "29.893.46"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.895.00"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7731"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3125"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.72"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07213"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08892"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.72%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.951.85"
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.430"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.188"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "29.937.60"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007855"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.155.97"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.163"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1593"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.518"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.044"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.428"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.05%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.559"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.08%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.544"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.109"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05495"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.248"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7503"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9985"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.712"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.64%  "
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4503"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").Value = "1.093.69"
$ws.Range("E43").Value = "  -4.66%  "
$ws.Range("E44").Value = "  +2.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8543"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.886"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.630"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.822"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.966"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.37%  "
